$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")
$ws.Columns("B").Delete()
